$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trade data (row 5) matching the existing table structure
$ws.Cells.Item(5, 1).Value = 9975.93
$ws.Cells.Item(5, 2).Value = 10017
$ws.Cells.Item(5, 3).Value = 80.11
$ws.Cells.Item(5, 4).Value = 79.78
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(5, 6).Value = -0.41
$ws.Cells.Item(5, 7).Value = 42609.505243055559
$ws.Cells.Item(5, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(5, 8).Value = $false
